$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the M_exp (column B) values with new numbers
$ws.Range("B2").Value = 31.2
$ws.Range("B3").Value = 35.4
$ws.Range("B4").Value = 20.4
$ws.Range("B5").Value = 44
$ws.Range("B6").Value = 19.8
$ws.Range("B7").Value = 46.4
$ws.Range("B8").Value = 21.6
$ws.Range("B9").Value = 8.8

# New column D header and asymmetric-error data
$ws.Range("D1").Value = "M_exp_err"
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 2.6
$ws.Range("D4").Value = 1.5
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 5.3
$ws.Range("D8").Value = 1.8
$ws.Range("D9").Value = 1

# Keep selection consistent with the saved file (last active cell D9)
$ws.Range("D9").Select()
